$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginData")
$ws2 = $wb.Worksheets.Item("Environment")

$ws1.Range("B1").Value = "Gaurav123"
$ws2.Range("A1").Value = "https://classic.crmpro.com/"

$ws2.Range("A5").Select()
$ws1.Range("H18").Select()
